# hotfix double-encoding-issue by using triple-braces
#
# The template placeholders like {{Field__c}} were being double-encoded
# downstream, so every merge-field token in the contract template gets
# wrapped in triple braces {{{Field__c}}} instead of double braces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sentence that introduces the parties to the contract (甲 = AccountName)
$ws.Range("B4").Value = "{{{AccountName__c}}}（以下、甲という。）と　株式会社サンプル（以下、乙という。）は、"

# 氏名 (name) value cell
$ws.Range("C7").Value = "{{{AccountName__c}}}"

# 現住所 (address) value cell
$ws.Range("C8").Value = "{{{AccountAddress__c}}}"

# 雇用期間 (employment period)
$ws.Range("C11").Value = "{{{StartDateFormat__c}}} 〜 {{{EndDateFormat__c}}}  "

# 勤務場所 (work location)
$ws.Range("C12").Value = "{{{Address__c}}}"

# 仕事内容 (job description)
$ws.Range("C13").Value = "{{{JobDescription__c}}} "

# 就業時間 (working hours)
$ws.Range("C14").Value = "{{{StartTime__c}}} 〜 {{{EndTime__c}}}  "

# 所定外労働の有無 (overtime y/n)
$ws.Range("C15").Value = "{{{hasOverTime__c}}}"

# 休暇 (holiday type)
$ws.Range("C16").Value = "{{{HoliDayType__c}}} "

# 賃金 (salary)
$ws.Range("C17").Value = "基本給(月)　{{{Salary__c}}}万円"

# 賃金締切日 (pay cutoff date)
$ws.Range("C18").Value = "{{{DueDate__c}}} "

# 賃金支払日 (pay date)
$ws.Range("C19").Value = "{{{SalaryDate__c}}} "

# Leave the same selection state the author had when they saved the file
$ws.Range("C13:H13").Select()
